# Clear the recorded "power (mW)" (column AC) readings for rows 2-11.
# These values are no longer available/valid for these spectra, so the
# cells become blank (matching the upstream change that dropped
# AC2:AC11 down to empty cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AC2:AC11").ClearContents()
